$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $row, $col, $val) {
    if ($null -eq $val -or $val -eq "") {
        $ws.Cells.Item($row, $col).Value = $null
    } else {
        $ws.Cells.Item($row, $col).Value = "'" + $val
    }
}

function Set-NumCell($ws, $row, $col, $val) {
    if ($null -eq $val) {
        $ws.Cells.Item($row, $col).Value = $null
    } else {
        $ws.Cells.Item($row, $col).Value = $val
    }
}

# ---- 1. Bulk-update the timestamp column (O2:O348) ----
for ($r = 2; $r -le 348; $r++) {
    $ws.Cells.Item($r, 15).Value = "2022-12-18 20:49:35"
}

# ---- 2. Row data for swapped / rotated product rows ----
$rowData = @{}
$rowData[18] = @("6444276", "Frischback Bio Buttergipfel", "/de/lebensmittel/brot-backwaren/haltbare-brote/frischback-brot/frischback-bio-buttergipfel/p/6444276", 12, 4, "Coop", "3.50", "2.03/100g", "Preis pro 100 Gramm", "2.03", "100g", "['lebensmittel', 'brot-backwaren', 'haltbare-brote', 'frischback-brot']", "Frischback Bio Buttergipfel - Online kein Bestand 3.50 Schweizer Franken", "['chilled']")
$rowData[19] = @("4913403", "Butterzopf", "/de/lebensmittel/brot-backwaren/baeckerei/butterzopf/p/4913403", 14, 4.5, "Coop", "2.80", "0.88/100g", "Preis pro 100 Gramm", "0.88", "100g", "['lebensmittel', 'brot-backwaren', 'baeckerei']", "Butterzopf 2.80 Schweizer Franken", $null)
$rowData[107] = @("4835222", "Reiswaffeln Barbecue Flavour", "/de/lebensmittel/brot-backwaren/cracker-knaeckebrot/reiswaffeln/reiswaffeln-barbecue-flavour/p/4835222", 29, 5, "Coop", "2.70", "3.38/100g", "Preis pro 100 Gramm", "3.38", "100g", "['lebensmittel', 'brot-backwaren', 'cracker-knaeckebrot', 'reiswaffeln']", "Reiswaffeln Barbecue Flavour 2.70 Schweizer Franken", "['gluten_free']")
$rowData[108] = @("4004924", "Maiswaffeln Quinoa &amp; Amaranth", "/de/lebensmittel/brot-backwaren/cracker-knaeckebrot/reiswaffeln/maiswaffeln-quinoa-amaranth/p/4004924", 35, 5, "Coop", "2.50", "1.92/100g", "Preis pro 100 Gramm", "1.92", "100g", "['lebensmittel', 'brot-backwaren', 'cracker-knaeckebrot', 'reiswaffeln']", "Maiswaffeln Quinoa &amp; Amaranth 2.50 Schweizer Franken", "['gluten_free']")
$rowData[126] = @("6825720", "Prix Garantie Roggenvollkornbrot", "/de/lebensmittel/brot-backwaren/haltbare-brote/lang-haltbare-brote/prix-garantie-roggenvollkornbrot/p/6825720", 5, 4, "Coop", "1.80", "0.36/100g", "Preis pro 100 Gramm", "0.36", "100g", "['lebensmittel', 'brot-backwaren', 'haltbare-brote', 'lang-haltbare-brote']", "Prix Garantie Roggenvollkornbrot 1.80 Schweizer Franken", $null)
$rowData[127] = @("4750751", "Roland Knäckebrot Hafer", "/de/lebensmittel/brot-backwaren/cracker-knaeckebrot/knaeckebrot-pancroc/roland-knaeckebrot-hafer/p/4750751", 19, 4.5, "Roland", "3.75", "1.63/100g", "Preis pro 100 Gramm", "1.63", "100g", "['lebensmittel', 'brot-backwaren', 'cracker-knaeckebrot', 'knaeckebrot-pancroc']", "Roland Knäckebrot Hafer 3.75 Schweizer Franken", "['vegan', 'vegetarian', 'gluten_free']")
$rowData[209] = @("3041815", "Hug Zwieback Original 250G", "/de/lebensmittel/brot-backwaren/cracker-knaeckebrot/zwieback/hug-zwieback-original-250g/p/3041815", 28, 4.5, "Hug", "4.20", "1.68/100g", "Preis pro 100 Gramm", "1.68", "100g", "['lebensmittel', 'brot-backwaren', 'cracker-knaeckebrot', 'zwieback']", "Hug Zwieback Original 250G 4.20 Schweizer Franken", $null)
$rowData[210] = @("6986883", "Betty Bossi Strudelteig", "/de/lebensmittel/brot-backwaren/teig-tortenboden/blaetterteig/betty-bossi-strudelteig/p/6986883", $null, 0, "Coop", "2.50", "2.08/100g", "Preis pro 100 Gramm", "2.08", "100g", "['lebensmittel', 'brot-backwaren', 'teig-tortenboden', 'blaetterteig']", "Betty Bossi Strudelteig 2.50 Schweizer Franken", "['chilled', 'vegan', 'vegetarian']")
$rowData[237] = @("6010377", "Cailler Schoggi Branchli Cake", "/de/lebensmittel/brot-backwaren/kuchen-suessgebaeck/kuchen-toertchen/cailler-schoggi-branchli-cake/p/6010377", 28, 3.5, "Cailler", "5.80", "1.45/100g", "Preis pro 100 Gramm", "1.45", "100g", "['lebensmittel', 'brot-backwaren', 'kuchen-suessgebaeck', 'kuchen-toertchen']", "Cailler Schoggi Branchli Cake 5.80 Schweizer Franken", $null)
$rowData[238] = @("5777498", "Old el Paso Soft Taco Shells Mini", "/de/lebensmittel/brot-backwaren/haltbare-brote/fladenbrote-tortillas/old-el-paso-soft-taco-shells-mini/p/5777498", 14, 4.5, "Old el Paso", "5.95", "4.10/100g", "Preis pro 100 Gramm", "4.10", "100g", "['lebensmittel', 'brot-backwaren', 'haltbare-brote', 'fladenbrote-tortillas']", "Old el Paso Soft Taco Shells Mini 5.95 Schweizer Franken", $null)
$rowData[250] = @("4967195", "Bonne Maman Financier Mandeln", "/de/lebensmittel/brot-backwaren/kuchen-suessgebaeck/kuchen-toertchen/bonne-maman-financier-mandeln/p/4967195", 4, 4.5, "Bonne Maman", "5.50", "3.14/100g", "Preis pro 100 Gramm", "3.14", "100g", "['lebensmittel', 'brot-backwaren', 'kuchen-suessgebaeck', 'kuchen-toertchen']", "Bonne Maman Financier Mandeln 5.50 Schweizer Franken", $null)
$rowData[251] = @("6313768", "Roland Petite Pause Choco Céréales", "/de/lebensmittel/suesses-snacks/guetzli-suessgebaeck/guetzli-mit-schokolade/schokolade-als-zutat/roland-petite-pause-choco-cereales/p/6313768", 2, 4.5, "Roland", "2.65", "2.52/100g", "Preis pro 100 Gramm", "2.52", "100g", "['lebensmittel', 'suesses-snacks', 'guetzli-suessgebaeck', 'guetzli-mit-schokolade', 'schokolade-als-zutat']", "Roland Petite Pause Choco Céréales 2.65 Schweizer Franken", $null)
$rowData[271] = @("4490655", "Brossard Zwieback Honig", "/de/lebensmittel/brot-backwaren/cracker-knaeckebrot/zwieback/brossard-zwieback-honig/p/4490655", 2, 4, "Brossard", "4.40", "1.26/100g", "Preis pro 100 Gramm", "1.26", "100g", "['lebensmittel', 'brot-backwaren', 'cracker-knaeckebrot', 'zwieback']", "Brossard Zwieback Honig 4.40 Schweizer Franken", $null)
$rowData[272] = @("4063481", "Lu Prince Cake &amp; Choc 5 Stück", "/de/lebensmittel/brot-backwaren/kuchen-suessgebaeck/kuchen-toertchen/lu-prince-cake-choc-5-stueck/p/4063481", 1, 4, "Lu", "3.95", "2.63/100g", "Preis pro 100 Gramm", "2.63", "100g", "['lebensmittel', 'brot-backwaren', 'kuchen-suessgebaeck', 'kuchen-toertchen']", "Lu Prince Cake &amp; Choc 5 Stück 3.95 Schweizer Franken", $null)
$rowData[273] = @("4144505", "Pandorino", "/de/lebensmittel/suesses-snacks/guetzli-suessgebaeck/suesse-broetchen-panettone/pandorino/p/4144505", 13, 5, "Coop", "1.60", "1.78/100g", "Preis pro 100 Gramm", "1.78", "100g", "['lebensmittel', 'suesses-snacks', 'guetzli-suessgebaeck', 'suesse-broetchen-panettone']", "Pandorino 20% pro 3 Aktion 1.60 Schweizer Franken", $null)
$rowData[292] = @("3845685", "Slow Food Roggenbrot", "/de/lebensmittel/brot-backwaren/haltbare-brote/lang-haltbare-brote/slow-food-roggenbrot/p/3845685", 39, 3.5, "Coop", "4.50", "0.90/100g", "Preis pro 100 Gramm", "0.90", "100g", "['lebensmittel', 'brot-backwaren', 'haltbare-brote', 'lang-haltbare-brote']", "Slow Food Roggenbrot 4.50 Schweizer Franken", $null)
$rowData[293] = @("6967944", "St. Michel Galettes Mou Choco Blanc", "/de/lebensmittel/brot-backwaren/kuchen-suessgebaeck/kuchen-toertchen/st-michel-galettes-mou-choco-blanc/p/6967944", 2, 4, "St Michel", "2.95", "1.64/100g", "Preis pro 100 Gramm", "1.64", "100g", "['lebensmittel', 'brot-backwaren', 'kuchen-suessgebaeck', 'kuchen-toertchen']", "St. Michel Galettes Mou Choco Blanc 2.95 Schweizer Franken", $null)

# ---- 3. Apply the row data: columns A-N, with correct text/number handling ----
# Column map (1-based): 1=A id(text) 2=B title(text) 3=C href(text) 4=D ratingAmount(num)
# 5=E ratingValue(num) 6=F brand(text) 7=G price(text) 8=H priceContext(text)
# 9=I priceContextHiddenText(text) 10=J priceContextPrice(text) 11=K priceContextAmount(text)
# 12=L udoCat(text) 13=M productAriaLabel(text) 14=N declarationIcons(text)
$textCols = @(1,2,3,6,7,8,9,10,11,12,13,14)
$numCols = @(4,5)

foreach ($row in $rowData.Keys) {
    $vals = $rowData[$row]
    foreach ($c in $textCols) {
        Set-TextCell $ws $row $c $vals[$c - 1]
    }
    foreach ($c in $numCols) {
        Set-NumCell $ws $row $c $vals[$c - 1]
    }
}

# ---- 4. Standalone M-column ("Online kein Bestand") text updates ----
$ws.Cells.Item(7, 13).Value = "'Betty Bossi Frischback Silserbuttergipfel IP-Suisse - Online kein Bestand 3.20 Schweizer Franken"
$ws.Cells.Item(34, 13).Value = "'Leisi Blätterteig rund ausgewallt Ø32cm - Online kein Bestand 3.20 Schweizer Franken"
$ws.Cells.Item(208, 13).Value = "'Betty Bossi Dinkel Blätterteig eckig - Online kein Bestand 3.50 Schweizer Franken"
